$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("A4").Value = '2026-01-22T10:28:38.550Z'
$ws.Range("B4").Value = 'John Doe'
$ws.Range("C4").Value = 'Male'
$ws.Range("D4").Value = 36
$ws.Range("E4").Value = 'Asthma Attack'
$ws.Range("F4").Value = '68e36fbc9950693f31349cc7'
$ws.Range("G4").Value = 'John Doe is currently experiencing an active asthma attack with several serious symptoms affecting breathing and daily activities. Vital signs remain stable, but symptoms suggest urgent care is needed.'
$ws.Range("H4").Value = 'HIGH: Severe difficulty breathing and inability to perform usual tasks detected. Immediate medical attention recommended. | MEDIUM: Increased coughing, faster breathing, and low peak flow noted, indicating worsening asthma symptoms.'
$ws.Range("I4").Value = 'Blood pressure is stable at 120/80 mmHg. | Heart rate remains steady at 70 bpm. | Oxygen saturation is within normal range at 98%. | Body temperature is normal at 98.6°F. | Blood sugar level is within normal limits at 100 mg/dL.'
$ws.Range("J4").Value = 'Morning meal routines show no significant impact on vitals. | Symptoms such as coughing and difficulty breathing increase throughout the day. | Reduced ability to perform usual tasks correlates with symptom severity.'
$ws.Range("K4").Value = 'Avoid exposure to known asthma triggers including smoke, dust, and allergens. | Maintain a calm environment to reduce breathing difficulties. | Consider monitoring peak flow regularly to detect changes early. | Stay hydrated and follow a balanced diet to support respiratory health.'
$ws.Range("L4").Value = 'Patient exhibits critical respiratory symptoms requiring close observation. | Vital signs are stable but symptoms indicate possible escalation. | Recommend reviewing current asthma management plan and medication adherence. | Emergency response may be necessary if symptoms worsen.'
$ws.Range("M4").Value = 'Schedule an urgent consultation with a healthcare provider. | Monitor symptom changes closely and seek immediate help if breathing worsens. | Keep emergency inhalers accessible at all times. | Plan for a follow-up appointment to reassess condition and treatment.'
$ws.Range("N4").Value = 1048
$ws.Range("O4").Value = 419
$ws.Range("P4").Value = 1467

# Row 5
$ws.Range("A5").Value = '2026-01-22T10:30:08.111Z'
$ws.Range("B5").Value = 'John Doe'
$ws.Range("C5").Value = 'Male'
$ws.Range("D5").Value = 36
$ws.Range("E5").Value = 'Asthma Attack'
$ws.Range("F5").Value = '68e36fbc9950693f31349cc7'
$ws.Range("G5").Value = 'John Doe is currently experiencing an active asthma attack with several severe symptoms. Vital signs are within normal ranges but symptoms indicate urgent careful monitoring.'
$ws.Range("H5").Value = 'HIGH: Severe asthma symptoms including inability to perform tasks, sleep, eat, and talk, along with faster breathing, require immediate attention.'
$ws.Range("I5").Value = 'Blood pressure stable at 120/80 mmHg. | Heart rate steady at 70 bpm. | Oxygen saturation normal at 98%. | Body temperature consistent at 98.6°F. | Blood sugar at 100 mg/dL is within normal range.'
$ws.Range("J5").Value = 'Routine daily activities like meals have been maintained but symptoms significantly impair usual functioning. | No significant fluctuations in vital signs observed during routine.'
$ws.Range("K5").Value = 'Avoid exposure to asthma triggers such as smoke, dust, or allergens. | Maintain a calm environment and try controlled breathing exercises to ease symptoms. | Ensure regular hydration and balanced nutrition despite current difficulties.'
$ws.Range("L5").Value = 'Patient is exhibiting high severity asthma symptoms. Careful observation and prompt follow-up are necessary. | Vital signs currently stable, but symptom severity suggests risk of deterioration. | Encourage patient to report any worsening symptoms immediately.'
$ws.Range("M5").Value = 'Schedule urgent consultation with respiratory specialist. | Consider in-person evaluation to assess lung function. | Monitor symptoms closely and seek emergency care if breathing worsens.'
$ws.Range("N5").Value = 1048
$ws.Range("O5").Value = 342
$ws.Range("P5").Value = 1390

# Row 6
$ws.Range("A6").Value = '2026-01-22T10:30:31.003Z'
$ws.Range("B6").Value = 'John Doe'
$ws.Range("C6").Value = 'Male'
$ws.Range("D6").Value = 36
$ws.Range("E6").Value = 'Asthma Attack'
$ws.Range("F6").Value = '68e36fbc9950693f31349cc7'
$ws.Range("G6").Value = 'John Doe is currently experiencing an active asthma attack with significant symptoms impacting breathing and daily activities. Vitals appear stable at the moment, but symptoms require close monitoring and immediate care.'
$ws.Range("H6").Value = 'HIGH: Severe asthma symptoms detected including faster breathing, inability to perform usual tasks, and difficulty breathing that need urgent medical attention.'
$ws.Range("I6").Value = 'Stable blood pressure at 120/80 mmHg. | Heart rate steady at 70 bpm. | Oxygen saturation healthy at 98%. | Normal body temperature at 98.6°F. | Blood sugar levels within normal range at 100 mg/dL.'
$ws.Range("J6").Value = 'Morning routine includes breakfast around 10:00 AM which may coincide with symptom onset. | Asthma symptoms such as coughing and breathing difficulty worsen throughout the day, affecting usual tasks.'
$ws.Range("K6").Value = 'Avoid exposure to asthma triggers like dust and smoke. | Practice slow, controlled breathing exercises to help manage symptoms. | Maintain a regular meal schedule avoiding heavy foods that may worsen breathing.'
$ws.Range("L6").Value = 'Symptom severity is high; urgent evaluation and possible adjustment of asthma management plan are recommended. | Patient should be advised on when to seek emergency care if symptoms worsen.'
$ws.Range("M6").Value = 'Schedule an urgent consultation with the pulmonologist or asthma specialist. | Monitor symptoms closely with daily logs for any changes. | Ensure access to prescribed inhalers and emergency asthma action plan.'
$ws.Range("N6").Value = 1048
$ws.Range("O6").Value = 350
$ws.Range("P6").Value = 1398

# Row 7
$ws.Range("A7").Value = '2026-01-22T10:34:23.486Z'
$ws.Range("B7").Value = 'John Doe'
$ws.Range("C7").Value = 'Male'
$ws.Range("D7").Value = 36
$ws.Range("E7").Value = 'Asthma Attack'
$ws.Range("F7").Value = '68e36fbc9950693f31349cc7'
$ws.Range("G7").Value = 'John Doe is currently experiencing an active asthma attack with multiple symptoms impacting his daily functions. Vital signs remain stable, but symptom severity indicates close monitoring is needed.'
$ws.Range("H7").Value = 'HIGH: Severe breathing difficulties and inability to perform usual tasks detected. Immediate attention recommended.'
$ws.Range("I7").Value = 'Blood pressure steady at 120/80 mmHg. | Heart rate stable at 70 bpm. | Oxygen saturation healthy at 98%. | Body temperature normal at 98.6°F. | Blood sugar level consistent at 100 mg/dL.'
$ws.Range("J7").Value = 'Symptom severity increases throughout the day affecting sleep and daily activities. | Routine meals are maintained but breathing difficulties may hinder eating comfort.'
$ws.Range("K7").Value = 'Avoid exposure to asthma triggers and maintain indoor air quality. | Practice controlled breathing exercises to ease respiratory distress. | Maintain regular meal times and stay hydrated. | Ensure a restful environment to improve sleep quality.'
$ws.Range("L7").Value = 'Patient is showing significant symptoms consistent with an asthma attack. No changes in vitals but symptom severity requires vigilant observation. | Recommend evaluation by respiratory specialist for advanced management and treatment plan review.'
$ws.Range("M7").Value = 'Schedule an urgent consultation with the healthcare provider. | Monitor symptoms closely and seek emergency care if breathing worsens. | Consider pulmonary function testing to assess lung status.'
$ws.Range("N7").Value = 1048
$ws.Range("O7").Value = 331
$ws.Range("P7").Value = 1379

# Row 8
$ws.Range("A8").Value = '2026-01-22T10:37:56.041Z'
$ws.Range("B8").Value = 'John Doe'
$ws.Range("C8").Value = 'Male'
$ws.Range("D8").Value = 36
$ws.Range("E8").Value = 'Asthma Attack'
$ws.Range("F8").Value = '68e36fbc9950693f31349cc7'
$ws.Range("G8").Value = 'John Doe is currently experiencing an active asthma attack with multiple severe symptoms impacting daily life. Vital signs are stable, but symptoms such as faster breathing and inability to perform usual tasks require close monitoring.'
$ws.Range("H8").Value = 'HIGH: Severe asthma attack symptoms detected including faster breathing, inability to perform tasks, and difficulty breathing. Immediate medical attention is advised.'
$ws.Range("I8").Value = 'Blood pressure is normal at 120/80 mmHg. | Heart rate is stable at 70 bpm. | Oxygen saturation is within normal range at 98%. | Body temperature is normal at 98.6°F. | Blood sugar is stable at 100 mg/dL.'
$ws.Range("J8").Value = 'Symptoms worsen during usual daily activities, affecting performance and sleep. | Breakfast routine is consistent but symptoms could impact appetite and energy.'
$ws.Range("K8").Value = 'Avoid exposure to asthma triggers such as smoke, allergens, and strong odors. | Maintain regular hydration and a balanced diet to support lung health. | Practice controlled breathing exercises to help ease shortness of breath. | Ensure a calm and clean environment to promote better sleep.'
$ws.Range("L8").Value = 'Patient shows significant respiratory distress and requires urgent review. | Recommend close supervision and possibly adjusting asthma management plan. | Monitor symptoms closely to prevent worsening and emergency situations.'
$ws.Range("M8").Value = 'Schedule an urgent consultation with the respiratory specialist. | Consider emergency care if symptoms such as inability to talk or eat worsen. | Regularly track symptoms and vital signs to update the care team.'
$ws.Range("N8").Value = 1048
$ws.Range("O8").Value = 367
$ws.Range("P8").Value = 1415

# Row 9
$ws.Range("A9").Value = '2026-01-22T10:39:34.068Z'
$ws.Range("B9").Value = 'John Doe'
$ws.Range("C9").Value = 'Male'
$ws.Range("D9").Value = 36
$ws.Range("E9").Value = 'Asthma Attack'
$ws.Range("F9").Value = '68e36fbc9950693f31349cc7'
$ws.Range("G9").Value = 'John Doe is currently experiencing an active asthma attack with multiple severe symptoms affecting his breathing and daily functions. Vital signs remain stable at this time. Close monitoring and timely care are essential.'
$ws.Range("H9").Value = 'HIGH: Severe symptoms including inability to perform usual tasks, inability to sleep, and faster breathing require urgent attention.'
$ws.Range("I9").Value = 'Blood pressure is normal at 120/80 mmHg. | Heart rate steady at 70 bpm. | Oxygen saturation is adequate at 98%. | Body temperature is normal at 98.6°F. | Blood sugar level stable at 100 mg/dL.'
$ws.Range("J9").Value = 'Routine meals such as breakfast are maintained but asthma symptoms are impacting daily activities and sleep. | Increased symptom severity correlates with difficulty performing usual tasks and resting.'
$ws.Range("K9").Value = 'Avoid known asthma triggers and stay in a well-ventilated environment. | Maintain regular meal times and hydration to support overall health. | Incorporate gentle breathing exercises to help ease breathing discomfort.'
$ws.Range("L9").Value = 'Patient exhibits critical asthma symptoms despite stable vital signs. | Recommend continuous monitoring and preparation for escalation of care if symptoms worsen.'
$ws.Range("M9").Value = 'Schedule an urgent consultation with the healthcare provider for asthma management. | Consider emergency evaluation if breathing difficulties increase or oxygen levels drop. | Follow up regularly to adjust care plans based on symptom progression.'
$ws.Range("N9").Value = 1048
$ws.Range("O9").Value = 342
$ws.Range("P9").Value = 1390

# Row 10
$ws.Range("A10").Value = '2026-01-22T10:40:07.590Z'
$ws.Range("B10").Value = 'John Doe'
$ws.Range("C10").Value = 'Male'
$ws.Range("D10").Value = 36
$ws.Range("E10").Value = 'Asthma Attack'
$ws.Range("F10").Value = '68e36fbc9950693f31349cc7'
$ws.Range("G10").Value = 'John Doe is currently experiencing an active asthma attack with multiple severe symptoms impacting daily function. Vital signs remain stable, but symptoms suggest urgent care needs.'
$ws.Range("H10").Value = 'HIGH: Severe asthma symptoms detected including difficulty breathing, inability to perform usual tasks, and faster breathing. Urgent medical attention is advised.'
$ws.Range("I10").Value = 'Blood pressure is normal at 120/80 mmHg. | Heart rate steady at 70 bpm. | Oxygen saturation is within normal range at 98%. | Body temperature is normal at 98.6°F. | Blood sugar levels are stable at 100 mg/dL.'
$ws.Range("J10").Value = 'Symptom severity appears to disrupt sleep and daily activities significantly. | Routine meals such as breakfast are maintained, but ability to eat is affected by symptoms.'
$ws.Range("K10").Value = 'Avoid exposure to asthma triggers such as dust, smoke, and allergens. | Practice controlled breathing techniques to help reduce breathlessness. | Ensure a calm and restful environment to improve sleep quality. | Stay hydrated and maintain balanced nutrition with easy-to-eat meals.'
$ws.Range("L10").Value = 'Patient shows critical signs of an asthma attack requiring careful monitoring. | Vital signs are stable but symptom scores indicate need for immediate evaluation. | Recommend review of current asthma management plan and emergency response readiness.'
$ws.Range("M10").Value = 'Seek urgent consultation with a healthcare provider to assess current episode. | Consider scheduling a follow-up appointment for asthma management adjustment. | Prepare an emergency action plan including medication review and trigger avoidance.'
$ws.Range("N10").Value = 1048
$ws.Range("O10").Value = 368
$ws.Range("P10").Value = 1416

# Row 11
$ws.Range("A11").Value = '2026-01-22T10:40:43.853Z'
$ws.Range("B11").Value = 'John Doe'
$ws.Range("C11").Value = 'Male'
$ws.Range("D11").Value = 36
$ws.Range("E11").Value = 'Asthma Attack'
$ws.Range("F11").Value = '68e36fbc9950693f31349cc7'
$ws.Range("G11").Value = 'John Doe is currently experiencing an active asthma attack with multiple severe symptoms affecting his breathing and daily functions. Vital signs remain stable at this time, but symptoms indicate urgent attention is important to manage the condition.'
$ws.Range("H11").Value = 'HIGH: Severe asthma symptoms detected including inability to perform usual tasks, inability to sleep, and faster breathing. Immediate monitoring and care are advised.'
$ws.Range("I11").Value = 'Vital signs such as blood pressure (120/80), heart rate (70 bpm), oxygen level (98%), body temperature (98.6°F), and blood sugar (100 mg/dL) are within normal range. | Stable vitals currently do not reflect the severity of symptoms but continuous monitoring is essential.'
$ws.Range("J11").Value = 'Routine activities such as meals are maintained but symptoms like increased coughing and difficulty breathing may interfere with daily functions and sleep quality.'
$ws.Range("K11").Value = 'Avoid exposure to asthma triggers such as smoke, dust, and allergens. | Maintain a calm environment to prevent exacerbation of breathing difficulties. | Ensure timely and regular intake of prescribed asthma control medications as directed by a healthcare provider. | Incorporate gentle breathing exercises to help improve lung function if tolerated.'
$ws.Range("L11").Value = 'The care team is aware of the active asthma attack and is monitoring symptoms closely. | Patient should present to urgent care if shortness of breath worsens or if unable to speak or eat. | No changes in vital signs currently, but symptom severity suggests need for careful observation.'
$ws.Range("M11").Value = 'Schedule a follow-up consultation with a pulmonologist or asthma specialist as soon as possible. | Consider immediate evaluation at emergency care if symptoms such as inability to talk or eat persist or worsen. | Maintain a symptom diary to track changes and report to care team regularly.'
$ws.Range("N11").Value = 1048
$ws.Range("O11").Value = 413
$ws.Range("P11").Value = 1461

# Row 12
$ws.Range("A12").Value = '2026-01-22T10:46:10.444Z'
$ws.Range("B12").Value = 'John Doe'
$ws.Range("C12").Value = 'Male'
$ws.Range("D12").Value = 36
$ws.Range("E12").Value = 'Asthma Attack'
$ws.Range("F12").Value = '68e36fbc9950693f31349cc7'
$ws.Range("G12").Value = 'John Doe is currently experiencing an active asthma attack with multiple severe symptoms impacting daily activities and comfort. Vital signs remain stable for now, but symptoms indicate urgent need for care.'
$ws.Range("H12").Value = 'HIGH: Severe asthma symptoms including faster breathing, inability to perform usual tasks, and difficulty breathing detected. Immediate medical attention recommended.'
$ws.Range("I12").Value = 'Blood pressure steady at 120/80 mmHg. | Heart rate normal at 70 bpm. | Oxygen saturation healthy at 98%. | Body temperature consistent at 98.6°F. | Blood sugar level within normal range at 100 mg/dL.'
$ws.Range("J12").Value = 'Patient attended breakfast at usual time without reported symptoms related to eating. | Routine activities appear limited by severity of symptoms, leading to reduced ability to perform usual tasks.'
$ws.Range("K12").Value = 'Avoid known asthma triggers and maintain access to inhalers or emergency medication. | Ensure adequate rest and avoid strenuous activity during symptom flare-ups. | Maintain a clean, allergen-free environment to reduce asthma attacks. | Stay hydrated and monitor symptoms closely.'
$ws.Range("L12").Value = 'Patient shows signs of significant respiratory distress. Monitoring and timely intervention are crucial. | Encourage emergency consultation if symptoms worsen or do not improve promptly.'
$ws.Range("M12").Value = 'Schedule urgent consultation with pulmonologist or emergency care provider. | Monitor breathing patterns and symptom severity throughout the day. | Prepare for possible immediate treatment based on symptom progression.'
$ws.Range("N12").Value = 1048
$ws.Range("O12").Value = 352
$ws.Range("P12").Value = 1400

# Row 13
$ws.Range("A13").Value = '2026-01-22T11:00:41.821Z'
$ws.Range("B13").Value = 'John Doe'
$ws.Range("C13").Value = 'Male'
$ws.Range("D13").Value = 36
$ws.Range("E13").Value = 'Asthma Attack'
$ws.Range("F13").Value = '68e36fbc9950693f31349cc7'
$ws.Range("G13").Value = 'John Doe is currently experiencing an active asthma attack with several severe symptoms affecting breathing and daily activities. Vital signs are stable at present, but symptoms indicate urgent attention is needed. This overview is informational and not a medical diagnosis.'
$ws.Range("H13").Value = 'HIGH: Severe asthma symptoms detected including inability to perform usual tasks, inability to eat or sleep, and faster breathing. | MEDIUM: Increased coughing, difficulty breathing, and low peak flow noted; monitor closely and seek care if symptoms worsen.'
$ws.Range("I13").Value = 'Vital signs such as blood pressure, heart rate, oxygen level, body temperature, and blood sugar are within normal ranges currently. | No abnormal vital sign fluctuations detected during recent assessments.'
$ws.Range("J13").Value = 'Symptoms worsen as daily tasks become difficult and sleep and eating are impaired. | Routine activities such as meals appear consistent, but respiratory symptoms impact overall well-being.'
$ws.Range("K13").Value = 'Avoid known asthma triggers such as allergens and smoke. | Maintain regular breathing exercises to support lung function. | Ensure consistent medication adherence as prescribed by your healthcare provider. | Keep hydrated and rest as symptoms permit.'
$ws.Range("L13").Value = 'Patient shows critical asthma symptoms that require close monitoring. | Encourage patient to report any worsening breathing difficulty immediately. | Support patient in managing lifestyle factors that impact asthma control.'
$ws.Range("M13").Value = 'Schedule an urgent consultation with your healthcare provider to reassess asthma management. | Consider visiting emergency care if symptoms like inability to talk or severe breathlessness persist or worsen. | Monitor symptoms and vitals regularly and maintain communication with your care team.'
$ws.Range("N13").Value = 1048
$ws.Range("O13").Value = 389
$ws.Range("P13").Value = 1437
